$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the existing hyperlink collection up front. Row insertion shifts
# the <hyperlink ref="..."> cell addresses but not which rId (i.e. which
# URL) belongs to which row, so links are rebuilt from scratch afterwards.
$ws.Hyperlinks.Delete()

# Insert a new row above the old row 4 (Stable Diffusion posting) to make
# room for the new PoC listing; rows 4-10 shift down to rows 5-11.
$ws.Rows.Item(4).Insert()

# Rewrite every data row (2-11) explicitly so cell values are correct
# regardless of how the insert shifted the previous contents.
$ws.Range("A2").Value = '2025-11-17 18:25:55'
$ws.Range("B2").Value = '【Next.js × TypeScript × Tailwind】コンポーネント制作パートナー募集!'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5428507'
$ws.Range("G2").Value = 528
$ws.Range("H2").Value = '🔥AI,Next.js'

$ws.Range("A3").Value = '2025-11-17 18:25:55'
$ws.Range("B3").Value = '大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5423720'
$ws.Range("G3").Value = 385
$ws.Range("H3").Value = '🔥AI,Ai ◆効率化'

$ws.Range("A4").Value = '2025-11-17 18:25:55'
$ws.Range("B4").Value = '【AI分析】企業利益最大化のためのPoC開発依頼'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5434977'
$ws.Range("G4").Value = 375
$ws.Range("H4").Value = '🔥AI,Ai ◆開発'

$ws.Range("A5").Value = '2025-11-17 18:25:55'
$ws.Range("B5").Value = 'Stable Diffusionに詳しいLoRAなどを用いた画像生成AIエンジニア募集'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5416328'
$ws.Range("G5").Value = 310
$ws.Range("H5").Value = '🔥AI,Ai'

$ws.Range("A6").Value = '2025-11-17 18:25:55'
$ws.Range("B6").Value = '製造業向けAI戦略アドバイザー募集(事業価値試算・プロジェクト推進支援)'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5419380'
$ws.Range("G6").Value = 298
$ws.Range("H6").Value = '🔥AI,Ai'

$ws.Range("A7").Value = '2025-11-17 18:25:55'
$ws.Range("B7").Value = '医療系機械学習モデル活用のGUIアプリ開発'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5435875'
$ws.Range("G7").Value = 225
$ws.Range("H7").Value = '🔥機械学習 ◆開発 ◇アプリ'

$ws.Range("A8").Value = '2025-11-17 18:25:55'
$ws.Range("B8").Value = '【自動運転プロジェクト経験者募集】実証実験・開発を推進するプロジェクトマネージャー'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5431107'
$ws.Range("G8").Value = 68
$ws.Range("H8").Value = '◆開発'

$ws.Range("A9").Value = '2025-11-17 18:25:55'
$ws.Range("B9").Value = 'UTAGE構築代行|ヒアリングから構築までお任せしたいです。'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5429882'
$ws.Range("G9").Value = 18
$ws.Range("H9").ClearContents()

$ws.Range("A10").Value = '2025-11-17 18:25:55'
$ws.Range("B10").Value = '【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '20,000 円 ~ 30,000 円 / 募集期間 5 日、取引期間 0 日'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5435667'
$ws.Range("G10").Value = 10
$ws.Range("H10").ClearContents()

$ws.Range("A11").Value = '2025-11-17 18:25:55'
$ws.Range("B11").Value = 'ロリポップ!レンタルサーバーの不具合を解決したい'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5435519'
$ws.Range("G11").Value = 10
$ws.Range("H11").ClearContents()

# Re-create the hyperlinks on column F in row order so relationship IDs
# (rId1..rId10) line up with F2..F11 the same way the original file did.
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5428507') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5423720') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5434977') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5416328') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5419380') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5435875') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5431107') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5429882') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5435667') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5435519') | Out-Null

# Hyperlinks.Add() stamps its own ad-hoc "Hyperlink font" style variant;
# reapply the workbook's built-in Hyperlink cell style so F2:F11 keep the
# same single shared style the sheet used before (and after) the edit.
$ws.Range("F2:F11").Style = "Hyperlink"

